$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 currently holds the blank/masked placeholder record (personnel_id 0,
# empty text fields, masked email/password) and row 11 holds the "Test"
# record (personnel_id 10). This review swaps their roles: row 2 becomes the
# real "Max Mustermann" sample user, and row 11 becomes the blank/masked
# placeholder that row 2 used to be.

# Write the real data into row 2.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "max_mustermann"
$ws.Range("C2").Value = "Max"
$ws.Range("D2").Value = "Mustermann"
$ws.Range("E2").Value = "Musterstraße 1"
$ws.Range("F2").Value = 12345
$ws.Range("G2").Value = "Musterstadt"
$ws.Range("H2").Value = "max.mustermann@example.com"
$ws.Range("I2").Value = "passwort123"
$ws.Range("J2").Value = 1

# Finish row 11: personnel_id/zip back to 0, text fields blanked out (kept as
# empty *text* values, not fully cleared cells, to match the blank-placeholder
# row's original cell typing), masked email/password, role_id 1.
$ws.Range("A11").Value = 0
$ws.Range("B11").Value = "'"
$ws.Range("C11").Value = "'"
$ws.Range("D11").Value = "'"
$ws.Range("E11").Value = "'"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "'"
$ws.Range("H11").Value = "..._...@...."
$ws.Range("I11").Value = "**********"
$ws.Range("J11").Value = 1
